# Update the EPEX Spot price workbook with the latest day of data.
#  - "Prix Spot" sheet: add a new date column AC ("12-jul") with 24 hourly prices.
#  - "Gaz" sheet: append a new row (26) for 2025-07-10.
#  - "CO2" sheet: append a new row (26) for 2025-07-10.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column AC, header "12-jul", rows 2..25 hourly values
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("AB1").Copy()
$wsPrix.Range("AC1").PasteSpecial(-4122)
$wsPrix.Range("AC1").Value = "12-jul"

$prixValues = @(100, 90.03, 82.67, 70.05, 63.95, 71.87, 69.45, 70.89, 65.59999999999999, 59.78, 18.19, 3, 0.22, 0, 0, 0.22, 13.68, 35, 50.08, 75.40000000000001, 107.75, 99.72, 117.17, 105.98)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 29).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 26 for 2025-07-10
# ---------------------------------------------------------------------------
# The date column stores plain text (e.g. "2025-07-09"), not a real Excel
# date serial. Assigning a YYYY-MM-DD-looking string directly to .Value
# gets auto-parsed into a date serial + date number format, so force a
# text format first and restore the default style afterwards to match the
# existing (unstyled) cells in the column.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A26").NumberFormat = "@"
$wsGaz.Range("A26").Value = "2025-07-10"
$wsGaz.Range("A26").Style = "Normal"
$wsGaz.Range("B26").Value = 34.4

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 26 for 2025-07-10
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A26").NumberFormat = "@"
$wsCo2.Range("A26").Value = "2025-07-10"
$wsCo2.Range("A26").Style = "Normal"
$wsCo2.Range("B26").Value = 69.8
